# Event template header cleanup:
#  - strip the trailing "*" from the "required field" header labels
#  - make the whole header row (row 1) bold
#  - move the active selection to L2 (matches the saved workbook's view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing "*" marker from the required-field headers in row 1.
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "location"
$ws.Range("D1").Value = "address"
$ws.Range("E1").Value = "city"
$ws.Range("H1").Value = "capacity"
$ws.Range("I1").Value = "startDate"
$ws.Range("J1").Value = "endDate"
$ws.Range("K1").Value = "type"
$ws.Range("L1").Value = "day1"

# Bold the entire header row (A1:U1).
$ws.Range("A1:U1").Font.Bold = $true

# Restore the saved selection/active cell.
$ws.Range("L2").Select() | Out-Null
